# Auto-generated Excel COM-interop edit script
# Applies cached market-data refresh values to the Leve profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
$wb = $excel.ActiveWorkbook

# --- ALC!row33 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 109.045456
$ws.Range("I33").Value = 117.333336
$ws.Range("K33").Value = 117.333336
$ws.Range("M33").Value = 111.666664

# --- ALC!row80 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 1287
$ws.Range("I80").Value = 1149.375
$ws.Range("J80").Value = 1837.5
$ws.Range("K80").Value = 3448.125
$ws.Range("L80").Value = 5512.5
$ws.Range("M80").Value = -2450.125
$ws.Range("N80").Value = -7508.5

# --- ALC!row83 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 1287
$ws.Range("I83").Value = 1149.375
$ws.Range("J83").Value = 1837.5
$ws.Range("K83").Value = 10344.375
$ws.Range("L83").Value = 16537.5
$ws.Range("M83").Value = -5352.375
$ws.Range("N83").Value = -26521.5

# --- ALC!row92 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 332.5
$ws.Range("I92").Value = 319.69232
$ws.Range("J92").Value = 499
$ws.Range("K92").Value = 319.69232
$ws.Range("L92").Value = 499
$ws.Range("M92").Value = 928.30768
$ws.Range("N92").Value = -2995

# --- ALC!row132 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1325.5714
$ws.Range("I132").Value = 1341.7273
$ws.Range("K132").Value = 4025.1819
$ws.Range("M132").Value = -1495.1819

# --- ALC!row138 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3045.182
$ws.Range("J138").Value = 4000
$ws.Range("L138").Value = 12000
$ws.Range("N138").Value = -22280

# --- ARM!row2 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2739.4546
$ws.Range("I2").Value = 1068.5
$ws.Range("J2").Value = 4744.6
$ws.Range("K2").Value = 1068.5
$ws.Range("L2").Value = 4744.6
$ws.Range("M2").Value = -955.5
$ws.Range("N2").Value = -4970.6

# --- ARM!row32 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2525.8286
$ws.Range("I32").Value = 1089.6428
$ws.Range("J32").Value = 8270.571
$ws.Range("K32").Value = 1089.6428
$ws.Range("L32").Value = 8270.571
$ws.Range("M32").Value = -802.6428000000001
$ws.Range("N32").Value = -8844.571

# --- ARM!row61 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2752.125
$ws.Range("I61").Value = 2388.5
$ws.Range("J61").Value = 3358.1667
$ws.Range("K61").Value = 2388.5
$ws.Range("L61").Value = 3358.1667
$ws.Range("M61").Value = -2176.5
$ws.Range("N61").Value = -3782.1667

# --- ARM!row63 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3566.6667
$ws.Range("I63").Value = 2880
$ws.Range("K63").Value = 2880
$ws.Range("M63").Value = -2194

# --- ARM!row66 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 3566.6667
$ws.Range("I66").Value = 2880
$ws.Range("K66").Value = 14400
$ws.Range("M66").Value = -10968

# --- ARM!row74 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2284.8572
$ws.Range("I74").Value = 2199
$ws.Range("J74").Value = 2499.5
$ws.Range("K74").Value = 2199
$ws.Range("L74").Value = 2499.5
$ws.Range("M74").Value = -1325
$ws.Range("N74").Value = -4247.5

# --- ARM!row77 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2284.8572
$ws.Range("I77").Value = 2199
$ws.Range("J77").Value = 2499.5
$ws.Range("K77").Value = 10995
$ws.Range("L77").Value = 12497.5
$ws.Range("M77").Value = -6627
$ws.Range("N77").Value = -21233.5

# --- ARM!row109 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").ClearContents()
$ws.Range("N109").Value = 0

# --- ARM!row116 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 2739.4546
$ws.Range("I116").Value = 1068.5
$ws.Range("J116").Value = 4744.6
$ws.Range("K116").Value = 1068.5
$ws.Range("L116").Value = 4744.6
$ws.Range("M116").Value = 1225.5
$ws.Range("N116").Value = -9332.6

# --- ARM!row136 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2752.125
$ws.Range("I136").Value = 2388.5
$ws.Range("J136").Value = 3358.1667
$ws.Range("K136").Value = 7165.5
$ws.Range("L136").Value = 10074.5001
$ws.Range("M136").Value = -4615.5
$ws.Range("N136").Value = -15174.5001

# --- BSM!row3 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2739.4546
$ws.Range("I3").Value = 1068.5
$ws.Range("J3").Value = 4744.6
$ws.Range("K3").Value = 1068.5
$ws.Range("L3").Value = 4744.6
$ws.Range("M3").Value = -954.5
$ws.Range("N3").Value = -4972.6

# --- BSM!row80 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 176.1
$ws.Range("I80").Value = 73
$ws.Range("J80").Value = 244.83333
$ws.Range("K80").Value = 73
$ws.Range("L80").Value = 244.83333
$ws.Range("M80").Value = 925
$ws.Range("N80").Value = -2240.83333

# --- BSM!row83 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H83").Value = 176.1
$ws.Range("I83").Value = 73
$ws.Range("J83").Value = 244.83333
$ws.Range("K83").Value = 365
$ws.Range("L83").Value = 1224.16665
$ws.Range("M83").Value = 4627
$ws.Range("N83").Value = -11208.16665

# --- BSM!row105 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2465.5
$ws.Range("I105").Value = 2379.2856
$ws.Range("J105").Value = 2666.6667
$ws.Range("K105").Value = 2379.2856
$ws.Range("L105").Value = 2666.6667
$ws.Range("M105").Value = -632.2856000000002
$ws.Range("N105").Value = -6160.6667

# --- CRP!row31 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1187.9
$ws.Range("J31").Value = 964.3333
$ws.Range("L31").Value = 964.3333
$ws.Range("N31").Value = -1554.3333

# --- CRP!row34 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1187.9
$ws.Range("J34").Value = 964.3333
$ws.Range("L34").Value = 964.3333
$ws.Range("N34").Value = -1368.3333

# --- CRP!row58 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1660.8572
$ws.Range("J58").Value = 1658.8
$ws.Range("L58").Value = 1658.8
$ws.Range("N58").Value = -2064.8

# --- CRP!row105 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 895.75
$ws.Range("I105").Value = 895.75
$ws.Range("K105").Value = 895.75
$ws.Range("M105").Value = 851.25

# --- CRP!row107 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 405.66666
$ws.Range("I107").Value = 231.8
$ws.Range("J107").Value = 623
$ws.Range("K107").Value = 231.8
$ws.Range("L107").Value = 623
$ws.Range("M107").Value = 1688.2
$ws.Range("N107").Value = -4463

# --- CRP!row134 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 3260.7
$ws.Range("I134").Value = 2961
$ws.Range("K134").Value = 8883
$ws.Range("M134").Value = -6348

# --- CRP!row136 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1660.8572
$ws.Range("J136").Value = 1658.8
$ws.Range("L136").Value = 4976.4
$ws.Range("N136").Value = -10076.4

# --- CUL!row2 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 188.7619
$ws.Range("I2").Value = 107
$ws.Range("J2").Value = 263.0909
$ws.Range("K2").Value = 642
$ws.Range("L2").Value = 1578.5454
$ws.Range("M2").Value = -529
$ws.Range("N2").Value = -1804.5454

# --- CUL!row80 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 6604.375
$ws.Range("J80").Value = 6467
$ws.Range("L80").Value = 19401
$ws.Range("N80").Value = -21273

# --- CUL!row83 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H83").Value = 6604.375
$ws.Range("J83").Value = 6467
$ws.Range("L83").Value = 58203
$ws.Range("N83").Value = -67563

# --- CUL!row132 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1450
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 1450
$ws.Range("K132").Value = 0
$ws.Range("L132").ClearContents()
$ws.Range("M132").Value = 13050
$ws.Range("N132").Value = -18110

# --- GSM!row80 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3006
$ws.Range("I80").Value = 3006
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 3006
$ws.Range("L80").ClearContents()
$ws.Range("M80").Value = -2008
$ws.Range("N80").Value = 0

# --- GSM!row83 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 3006
$ws.Range("I83").Value = 3006
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 15030
$ws.Range("L83").ClearContents()
$ws.Range("M83").Value = -10038
$ws.Range("N83").Value = 0

# --- GSM!row126 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 8284.4
$ws.Range("I126").Value = 8284.4
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 24853.2
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -22383.2

# --- LTW!row22 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2458.1667
$ws.Range("I22").Value = 2856.8572
$ws.Range("J22").Value = 1900
$ws.Range("K22").Value = 2856.8572
$ws.Range("L22").Value = 1900
$ws.Range("M22").Value = -2561.8572
$ws.Range("N22").Value = -2490

# --- LTW!row27 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 2458.1667
$ws.Range("I27").Value = 2856.8572
$ws.Range("J27").Value = 1900
$ws.Range("K27").Value = 2856.8572
$ws.Range("L27").Value = 1900
$ws.Range("M27").Value = -2749.8572
$ws.Range("N27").Value = -2114

# --- LTW!row46 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2704.8286
$ws.Range("I46").Value = 2588.2354
$ws.Range("J46").Value = 2814.9443
$ws.Range("K46").Value = 2588.2354
$ws.Range("L46").Value = 2814.9443
$ws.Range("M46").Value = -2400.2354
$ws.Range("N46").Value = -3190.9443

# --- LTW!row82 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2149.5
$ws.Range("I82").Value = 2079.4
$ws.Range("J82").Value = 2500
$ws.Range("K82").Value = 2079.4
$ws.Range("L82").Value = 2500
$ws.Range("M82").Value = -1718.4
$ws.Range("N82").Value = -3222

# --- LTW!row85 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 2149.5
$ws.Range("I85").Value = 2079.4
$ws.Range("J85").Value = 2500
$ws.Range("K85").Value = 2079.4
$ws.Range("L85").Value = 2500
$ws.Range("M85").Value = -831.4000000000001
$ws.Range("N85").Value = -4996

# --- LTW!row135 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H135").Value = 50000
$ws.Range("J135").Value = 50000
$ws.Range("L135").Value = 50000
$ws.Range("N135").Value = -60140

# --- WVR!row100 ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1276.5454
$ws.Range("I100").Value = 1255.25
$ws.Range("K100").Value = 2510.5
$ws.Range("M100").Value = -1969.5
